# Adds new gazetteer entries (new rows of English/Russian name & treasure
# pairs) to the "Лист1" worksheet, mirroring the author's upload of
# additional rows under the existing География/Сокровища/Имена tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Геогрфия block (columns I/J): fill in three rows that previously had
#     only L/M (Сокровища) and O/P (Имена) entries.
$ws.Range("I40").Value = "Marlfagh"
$ws.Range("J40").Value = "Марлфаг"
$ws.Range("I41").Value = "Tonintoul"
$ws.Range("J41").Value = "Тонинтул"
$ws.Range("I42").Value = "Dun Morbhaidh"
$ws.Range("J42").Value = "Дун-Морбхайд"

# --- Имена block (columns O/P): new names filled into previously-blank
#     cells within the already-shaded O42:P47 range.
$ws.Range("O45").Value = "Oswiu"
$ws.Range("P47").Value = "Тарайн"

# --- Имена block continues past the old end of the shaded range.
$ws.Range("O52").Value = "Onthloug"
$ws.Range("P52").Value = "Онтлуг"
$ws.Range("O53").Value = "Ciniath"
$ws.Range("P53").Value = "Киниат"

# --- Сокровища / Имена new rows 54-65.
$ws.Range("L54").Value = "Revenant Sword"
$ws.Range("M54").Value = "Не упокоенный мститель"
$ws.Range("O54").Value = "Seonaid"
$ws.Range("P54").Value = "Шинейд"

$ws.Range("L55").Value = "Liquid Empathy"
$ws.Range("M55").Value = "Жидкая чувствительность"
$ws.Range("O55").Value = "Marius"
$ws.Range("P55").Value = "Мариус"

$ws.Range("L56").Value = "Uzelin"
$ws.Range("M56").Value = "Узелин"

$ws.Range("L57").Value = "Demonface Coin"
$ws.Range("M57").Value = "Демоноликая монета"

$ws.Range("L58").Value = "Treeheart"
$ws.Range("M58").Value = "Сердце древа"

$ws.Range("L59").Value = "Throne Embers"
$ws.Range("M59").Value = "Угли трона"

$ws.Range("L60").Value = "Ghoul Drum"
$ws.Range("M60").Value = "Барабан гуля"

$ws.Range("L61").Value = "Wisdom Teeth"
$ws.Range("M61").Value = "Зубы мудрости"

$ws.Range("L62").Value = "Temporary Disintegration"
$ws.Range("M62").Value = "Временный распад"

$ws.Range("L63").Value = "Child of Stone"
$ws.Range("M63").Value = "Дитя камня"

$ws.Range("L64").Value = "The Blooded Board"
$ws.Range("M64").Value = "Окровавленная доска"

$ws.Range("L65").Value = "The Bound Hand"
$ws.Range("M65").Value = "Подневольная рука"

# M58:M63 picked up a direct "apply fill" style in the source (same visual
# as default, fillId=0) when the author re-saved the range.
$ws.Range("M58:M63").Interior.ColorIndex = 0

# P45 remains blank but is highlighted in yellow, matching the new swatch
# introduced alongside this block (new cellXfs entry: fillId=3, no font
# change).
$ws.Range("P45").Interior.Color = 65535

# Reflect the author's final view state: zoomed to 70%, scrolled so row 19
# is at the top, with Q45 as the active selection.
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("Q45").Select()
